# Apply the "arkusz_test" sheet update: replace row 2 contents and append
# rows 3-10 with new RMA / serial-number records (commit: "to wszystko ma byc").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arkusz_test")

# Columns: A=lp, B=RMA, C=Nazwa urzadzenia, D=Nr seryjny przyjety, E=Nr seryjny wydany
$rows = @(
    @(194, "2012221611-01", "SSDSC2KB480G701", "PHYS81960047480BGN", "PHYS744100Z4480BGN"),
    @(195, "2012221611-02", "SSDSC2KB480G701", "PHYS819300CW480BGN", "BTYS82110B73480BGN"),
    @(196, "2012221611-03", "SSDSC2KB480G701", "PHYS819300DB480BGN", "BTYS807101X8480BGN"),
    @(197, "2012221611-04", "SSDSC2KB480G701", "PHYS819301UT480BGN", "PHYS738002M6480BGN"),
    @(198, "2012221611-05", "SSDSC2KB480G701", "PHYS819301UW480BGN", "PHYS7375003K480BGN"),
    @(199, "2012221611-06", "SSDSC2KB480G701", "PHYS819300CZ480BGN", "BTYS82010L38480BGN"),
    @(200, "2012221611-07", "SSDSC2KB480G701", "PHYS819600HB480BGN", "PHYS738000HY480BGN"),
    @(201, "2012221611-08", "SSDSC2KB480G701", "PHYS816104DB480BGN", "PHYS7375009W480BGN"),
    @(202, "2012221611-09", "SSDSC2KB480G701", "PHYS819300NA480BGN", "PHYS738000HH480BGN")
)

$startRow = 2
$endRow = $startRow + $rows.Count - 1

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Column A (lp.) on the newly-added rows (3-10) should carry the same look as
# the existing A2 (bold, centered, thin-box border) - copy that formatting
# down without disturbing the values we just wrote.
$ws.Range("A2").Copy()
$ws.Range("A3:A" + $endRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column F stays present-but-empty (as it already is on row 2) for every new
# row; writing a bare "'" forces a real (empty) text cell instead of Excel
# deleting/ignoring an outright "" assignment, then ClearFormats drops the
# stray quote-prefix formatting so the cell matches row 2's plain empty cell.
$fRange = $ws.Range("F3:F" + $endRow)
$fRange.Value = "'"
$fRange.ClearFormats()
